$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells touched in this update, keyed "COL ROW" -> new value.
# Force Text number-format before writing so numeric-looking strings
# (e.g. "9.868", "1.000", "215.80") are stored verbatim instead of
# being auto-parsed into doubles (which would mangle precision /
# drop trailing zeros). Resetting Style to "Normal" afterwards drops
# the temporary Text format so the cell keeps the original (default)
# style index.

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '25.879.17'
Set-TextValue $ws.Range("E2") '  -0.56%  '

Set-TextValue $ws.Range("D3") '1.632.83'
Set-TextValue $ws.Range("E3") '  -0.59%  '

Set-TextValue $ws.Range("E4") '  -0.11%  '

Set-TextValue $ws.Range("D5") '215.80'
Set-TextValue $ws.Range("E5") '  +0.50%  '

Set-TextValue $ws.Range("D6") '0.5114'
Set-TextValue $ws.Range("E6") '  +0.35%  '

Set-TextValue $ws.Range("D7") '1.002'
Set-TextValue $ws.Range("E7") '  -0.08%  '

Set-TextValue $ws.Range("D8") '0.2563'
Set-TextValue $ws.Range("E8") '  -0.11%  '

Set-TextValue $ws.Range("D9") '0.06343'
Set-TextValue $ws.Range("E9") '  -0.22%  '

Set-TextValue $ws.Range("D10") '19.47'
Set-TextValue $ws.Range("E10") '  -0.44%  '

Set-TextValue $ws.Range("D11") '0.07786'
Set-TextValue $ws.Range("E11") '  +0.38%  '

Set-TextValue $ws.Range("D12") '4.248'
Set-TextValue $ws.Range("E12") '  -0.80%  '

Set-TextValue $ws.Range("D13") '1.634.40'
Set-TextValue $ws.Range("E13") '  -0.58%  '

Set-TextValue $ws.Range("D14") '1.858.37'
Set-TextValue $ws.Range("E14") '  -0.61%  '

Set-TextValue $ws.Range("D15") '0.5522'
Set-TextValue $ws.Range("E15") '  +1.57%  '

Set-TextValue $ws.Range("D16") '63.76'
Set-TextValue $ws.Range("E16") '  -0.89%  '

Set-TextValue $ws.Range("D17") '0.0₅7609'
Set-TextValue $ws.Range("E17") '  -1.39%  '

Set-TextValue $ws.Range("D18") '25.915.91'
Set-TextValue $ws.Range("E18") '  -0.53%  '

Set-TextValue $ws.Range("E19") '  -0.05%  '

Set-TextValue $ws.Range("B20") 'Uniswap'
Set-TextValue $ws.Range("C20") 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range("D20") '4.423'
Set-TextValue $ws.Range("E20") '  +0.11%  '

Set-TextValue $ws.Range("B21") 'BitcoinCash'
Set-TextValue $ws.Range("C21") 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Range("D21") '194.78'
Set-TextValue $ws.Range("E21") '  -1.20%  '

Set-TextValue $ws.Range("D22") '9.868'

Set-TextValue $ws.Range("D23") '6.027'
Set-TextValue $ws.Range("E23") '  +0.01%  '

Set-TextValue $ws.Range("E24") '  -0.16%  '

Set-TextValue $ws.Range("E25") '  +1.57%  '

Set-TextValue $ws.Range("D26") '142.15'

Set-TextValue $ws.Range("D27") '0.1255'
Set-TextValue $ws.Range("E27") '  +5.42%  '

Set-TextValue $ws.Range("D28") '6.763'
Set-TextValue $ws.Range("E28") '  -0.77%  '

Set-TextValue $ws.Range("D29") '15.54'
Set-TextValue $ws.Range("E29") '  -0.30%  '

Set-TextValue $ws.Range("E30") '  +0.66%  '

Set-TextValue $ws.Range("D31") '0.04914'
Set-TextValue $ws.Range("E31") '  +1.30%  '

Set-TextValue $ws.Range("E32") '  -0.30%  '

Set-TextValue $ws.Range("D33") '3.184'
Set-TextValue $ws.Range("E33") '  +0.67%  '

Set-TextValue $ws.Range("E34") '  +1.57%  '

Set-TextValue $ws.Range("D35") '2.375'
Set-TextValue $ws.Range("E35") '  +0.38%  '

Set-TextValue $ws.Range("D36") '0.8967'
Set-TextValue $ws.Range("E36") '  +0.07%  '

Set-TextValue $ws.Range("D37") '0.5522'
Set-TextValue $ws.Range("E37") '  +1.39%  '

Set-TextValue $ws.Range("D38") '2.544'
Set-TextValue $ws.Range("E38") '  -1.38%  '

Set-TextValue $ws.Range("D39") '1.116.61'
Set-TextValue $ws.Range("E39") '  -2.02%  '

Set-TextValue $ws.Range("D40") '0.01554'
Set-TextValue $ws.Range("E40") '  -0.52%  '

Set-TextValue $ws.Range("D41") '1.000'
Set-TextValue $ws.Range("E41") '  -0.19%  '

Set-TextValue $ws.Range("D42") '5.585'
Set-TextValue $ws.Range("E42") '  +3.55%  '

Set-TextValue $ws.Range("D43") '0.7951'
Set-TextValue $ws.Range("E43") '  -1.76%  '

Set-TextValue $ws.Range("D44") '97.78'
Set-TextValue $ws.Range("E44") '  -1.57%  '

Set-TextValue $ws.Range("D45") '1.769.02'
Set-TextValue $ws.Range("E45") '  -0.62%  '

Set-TextValue $ws.Range("D46") '0.0₈116'
Set-TextValue $ws.Range("E46") '  -10.05%  '

Set-TextValue $ws.Range("E47") '  -2.04%  '

Set-TextValue $ws.Range("D48") '1.003'
Set-TextValue $ws.Range("E48") '  +0.37%  '

Set-TextValue $ws.Range("D49") '54.83'
Set-TextValue $ws.Range("E49") '  +0.06%  '

Set-TextValue $ws.Range("D50") '0.05137'
Set-TextValue $ws.Range("E50") '  +1.54%  '

Set-TextValue $ws.Range("D51") '7.583'
Set-TextValue $ws.Range("E51") '  +3.40%  '
